# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) block for rows 16-60
# is flipped top-to-bottom: the period list (1607..2003) now runs in
# ascending order top->bottom (it used to run 2003..1607 descending), and
# the "Valor Mora" amounts follow the same row-for-row mirror (the value
# that used to sit 45-rows-from-the-top now sits in the mirrored row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 60

# Ascending period codes (YYMM), 1607 .. 2003 - the reverse of the
# previous (descending) order.
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$rowCount = $lastRow - $firstRow + 1

# Read the existing "Valor Mora" column top-to-bottom first, then write
# it back bottom-to-top (i.e. mirror the column) next to the new period
# labels above.
$valores = @()
for ($i = 0; $i -lt $rowCount; $i++) {
    $valores += $ws.Cells.Item($firstRow + $i, 6).Value2
}

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$rowCount - 1 - $i]
}
